# Auto-generated PowerShell Excel COM-interop script
# Applies the cryptos.xlsx data refresh described in the diff (GitHub Actions crypto price update)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '43.129.17'
$ws.Cells.Item(2, 5).Value = '  +5.10%  '

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '2.235.11'
$ws.Cells.Item(3, 5).Value = '  +2.92%  '

$ws.Cells.Item(4, 5).Value = '  +0.06%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '246.35'
$ws.Cells.Item(5, 5).Value = '  +3.96%  '

$ws.Cells.Item(6, 5).Value = '  +0.93%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '75.35'
$ws.Cells.Item(7, 5).Value = '  +7.65%  '

$ws.Cells.Item(8, 5).Value = '  -0.15%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.616'
$ws.Cells.Item(9, 5).Value = '  +6.53%  '

$ws.Cells.Item(10, 5).Value = '  +2.85%  '

$ws.Cells.Item(11, 5).Value = '  +0.63%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '55.53'
$ws.Cells.Item(12, 5).Value = '  +0.57%  '

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '6.96'
$ws.Cells.Item(13, 5).Value = '  +3.03%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '2.570.90'
$ws.Cells.Item(15, 5).Value = '  +2.91%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '14.74'
$ws.Cells.Item(16, 5).Value = '  +6.16%  '

$ws.Cells.Item(17, 2).Value = 'WrappedEther'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '2.246.23'
$ws.Cells.Item(17, 5).Value = '  +3.92%  '

$ws.Cells.Item(18, 2).Value = 'Polygon'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.814'
$ws.Cells.Item(18, 5).Value = '  +0.64%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '43.024.51'
$ws.Cells.Item(19, 5).Value = '  +4.95%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '0.0000105'
$ws.Cells.Item(20, 5).Value = '  +2.73%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '71.07'
$ws.Cells.Item(21, 5).Value = '  +0.98%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '5.98'
$ws.Cells.Item(22, 5).Value = '  +0.91%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '10.52'
$ws.Cells.Item(23, 5).Value = '  +5.30%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '230.56'
$ws.Cells.Item(24, 5).Value = '  +2.27%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.19'
$ws.Cells.Item(25, 5).Value = '  +11.22%  '

$ws.Cells.Item(26, 5).Value = '  -0.07%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '10.95'
$ws.Cells.Item(27, 5).Value = '  +0.04%  '

$ws.Cells.Item(28, 2).Value = 'Toncoin'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '2.31'
$ws.Cells.Item(28, 5).Value = '  +5.56%  '

$ws.Cells.Item(29, 2).Value = 'WEMIXToken'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '3.36'
$ws.Cells.Item(29, 5).Value = '  -5.30%  '

$ws.Cells.Item(30, 2).Value = 'PancakeSwap'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '2.25'
$ws.Cells.Item(30, 5).Value = '  +2.60%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '174.52'
$ws.Cells.Item(31, 5).Value = '  +4.95%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '36.95'
$ws.Cells.Item(32, 5).Value = '  +18.61%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '20.35'
$ws.Cells.Item(33, 5).Value = '  +2.34%  '

$ws.Cells.Item(34, 5).Value = '  +2.71%  '

$ws.Cells.Item(35, 5).Value = '  +4.10%  '

$ws.Cells.Item(36, 5).Value = '  +1.10%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.111'
$ws.Cells.Item(37, 5).Value = '  +6.26%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '4.36'
$ws.Cells.Item(38, 5).Value = '  +5.04%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.0332'
$ws.Cells.Item(39, 5).Value = '  +15.89%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '13.12'
$ws.Cells.Item(40, 5).Value = '  +5.13%  '

$ws.Cells.Item(41, 5).Value = '  +3.51%  '

$ws.Cells.Item(42, 5).Value = '  +3.41%  '

$ws.Cells.Item(43, 5).Value = '  +4.84%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '60.05'
$ws.Cells.Item(44, 5).Value = '  -0.33%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '105.64'
$ws.Cells.Item(45, 5).Value = '  +7.42%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '8.55'
$ws.Cells.Item(46, 5).Value = '  +3.29%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.0992'
$ws.Cells.Item(47, 5).Value = '  +2.11%  '

$ws.Cells.Item(48, 2).Value = 'ARBITRUM'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '1.11'
$ws.Cells.Item(48, 5).Value = '  +1.85%  '

$ws.Cells.Item(49, 2).Value = 'WOONetwork'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.442'
$ws.Cells.Item(49, 5).Value = '  +20.79%  '

$ws.Cells.Item(50, 5).Value = '  +3.62%  '

$ws.Cells.Item(51, 5).Value = '  +2.24%  '
